# liensMagasinPT.xlsx - "Add files via upload" edit
#
# The author re-saved the sheet after tidying up a handful of rows:
#  - several rows that were missing a (blank) B/D placeholder cell in the
#    "spans=1:4" block now get one, matching the style already used by all
#    the other blank placeholder cells in the sheet;
#  - a handful of cells that had picked up stray "debug" fonts (the two
#    Consolas fonts used nowhere else, and an inconsistent hyperlink font
#    on D28) get normalised back to the same styles used throughout the
#    rest of the sheet;
#  - row 30 had its "-BackStore-" label typed into column C instead of B
#    (compare with rows 12/13/21/31/36), so the label moves to B30 and C30
#    gets the usual "Rien pour le moment..." placeholder text;
#  - the selection left active in the sheet moved from D38 to E13.
#
# We reproduce all of this purely with native Range operations: blank
# "format only" cells are produced with Copy + PasteSpecial(xlPasteFormats)
# from a cell that already carries the right style (so we reuse existing
# style records instead of inventing new ones), and the row 30 fix is a
# couple of plain Value assignments.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) Rows that gained a new, still-empty placeholder cell (column B or D)
#    using the same plain style already used by every other blank
#    placeholder cell in the sheet (e.g. B6, B7, B9 ...).
# ---------------------------------------------------------------------
$blankCells = @("B8","D12","D13","B17","B18","B20","D21","B26","B27","B28","B29","D31","B33","B35","B40")
foreach ($addr in $blankCells) {
    $ws.Range("B6").Copy()
    $ws.Range($addr).PasteSpecial($xlPasteFormats)
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Cells that were using the stray Consolas "debug" fonts get put back
#    on the normal Arial body style used throughout the sheet.
# ---------------------------------------------------------------------
$toNormalStyle = @("D8","D17","D18","D20","D25","D26","D27","D29","D33","D35","D40","B19","C26","C27","C29","C35")
foreach ($addr in $toNormalStyle) {
    $ws.Range("C1").Copy()
    $ws.Range($addr).PasteSpecial($xlPasteFormats)
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3) D28 was the one hyperlink cell still on the old Aptos-Narrow
#    hyperlink font; align it with every other hyperlink cell (e.g. D2).
# ---------------------------------------------------------------------
$ws.Range("D2").Copy()
$ws.Range("D28").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4) B37 was left at the default style; give it the same body style as
#    its neighbours while keeping its existing text.
# ---------------------------------------------------------------------
$ws.Range("B6").Copy()
$ws.Range("B37").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 5) Row 30: the "-BackStore-" label had been typed into column C; move
#    it to column B (matching rows 12/13/21/31/36) and restore the usual
#    "Rien pour le moment..." text in C30.
# ---------------------------------------------------------------------
$ws.Range("B30").Value = "-BackStore-"
$ws.Range("B12").Copy()
$ws.Range("B30").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("C30").Value = "Rien pour le moment..."

# ---------------------------------------------------------------------
# 6) Leave the same cell selected as in the saved workbook.
# ---------------------------------------------------------------------
[void]$ws.Range("E13").Select()
